$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 427, shifting existing rows 427.. down by one.
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new market record.
$ws.Cells.Item(427, 1).Value  = 5
$ws.Cells.Item(427, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(427, 3).Value  = "Maule"
$ws.Cells.Item(427, 4).Value  = 45154
$ws.Cells.Item(427, 5).Value  = 7
$ws.Cells.Item(427, 6).Value  = 100112008
$ws.Cells.Item(427, 7).Value  = "Coliflor"
$ws.Cells.Item(427, 8).Value  = "Sin especificar"
$ws.Cells.Item(427, 9).Value  = "Primera"
$ws.Cells.Item(427, 10).Value = 4000
$ws.Cells.Item(427, 11).Value = 800
$ws.Cells.Item(427, 12).Value = 800
$ws.Cells.Item(427, 13).Value = 800
$ws.Cells.Item(427, 14).Value = "$/unidad"
$ws.Cells.Item(427, 15).Value = "Región del Maule"
$ws.Cells.Item(427, 16).Value = 800
$ws.Cells.Item(427, 17).Value = 1
$ws.Cells.Item(427, 18).Value = "Hortaliza"
